# Actualiza los horarios de la Línea 141 con el nuevo scrap (04:18:02)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Hoja "LP1912"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:18:02"
$ws1.Range("A3").Value = "Total filas: 14"

# Inserta dos filas en blanco antes de la fila 13 (la vieja fila 13 pasa a
# ser la fila 15, conservando sus datos sin cambios).
$ws1.Rows.Item(13).Insert()
$ws1.Rows.Item(13).Insert()

# Fila 13 (nueva)
$ws1.Cells.Item(13, 1).Value = "04:18:02"
$ws1.Cells.Item(13, 2).Value = "05:34"
$ws1.Cells.Item(13, 3).Value = "14_ABASTO"
$ws1.Cells.Item(13, 4).Value = 76
$ws1.Cells.Item(13, 5).Value = "LP1912"

# Fila 14 (nueva)
$ws1.Cells.Item(14, 1).Value = "04:18:02"
$ws1.Cells.Item(14, 2).Value = "05:35"
$ws1.Cells.Item(14, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(14, 4).Value = 77
$ws1.Cells.Item(14, 5).Value = "LP1912"

# Fila 16 (nueva, al final)
$ws1.Cells.Item(16, 1).Value = "04:18:02"
$ws1.Cells.Item(16, 2).Value = "05:46"
$ws1.Cells.Item(16, 3).Value = "15_ABASTO"
$ws1.Cells.Item(16, 4).Value = 88
$ws1.Cells.Item(16, 5).Value = "LP1912"

# Fila 17 (nueva, al final)
$ws1.Cells.Item(17, 1).Value = "04:18:02"
$ws1.Cells.Item(17, 2).Value = "06:05"
$ws1.Cells.Item(17, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(17, 4).Value = 107
$ws1.Cells.Item(17, 5).Value = "LP1912"

# Fila 18 (nueva, al final)
$ws1.Cells.Item(18, 1).Value = "04:18:02"
$ws1.Cells.Item(18, 2).Value = "06:12"
$ws1.Cells.Item(18, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(18, 4).Value = 114
$ws1.Cells.Item(18, 5).Value = "LP1912"

# Fila 19 (nueva, al final)
$ws1.Cells.Item(19, 1).Value = "04:18:02"
$ws1.Cells.Item(19, 2).Value = "06:14"
$ws1.Cells.Item(19, 3).Value = "225_HARAS DEL SUR"
$ws1.Cells.Item(19, 4).Value = 116
$ws1.Cells.Item(19, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Hoja "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:18:02"
$ws2.Range("A3").Value = "Total filas: 4"

# Fila 8 (nueva, al final)
$ws2.Cells.Item(8, 1).Value = "04:18:02"
$ws2.Cells.Item(8, 2).Value = "05:35"
$ws2.Cells.Item(8, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(8, 4).Value = 77
$ws2.Cells.Item(8, 5).Value = "LP1912"

# Fila 9 (nueva, al final)
$ws2.Cells.Item(9, 1).Value = "04:18:02"
$ws2.Cells.Item(9, 2).Value = "06:12"
$ws2.Cells.Item(9, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(9, 4).Value = 114
$ws2.Cells.Item(9, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Hoja "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:18:02"
$ws3.Range("A3").Value = "Total filas: 2"

# Fila 7 (nueva, al final)
$ws3.Cells.Item(7, 1).Value = "04:18:02"
$ws3.Cells.Item(7, 2).Value = "06:09"
$ws3.Cells.Item(7, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(7, 4).Value = 111
$ws3.Cells.Item(7, 5).Value = "L6173"
